# Backing up the most current
# Update the "Final Balance" (E) and "Cumulative Returns (%)" (F) columns
# in the backtest results sheet with refreshed backtest values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  E = 1306.472775194389;  F = 30.64727751943894 },
    @{ Row = 3;  E = 773.9697588173583;  F = -22.60302411826417 },
    @{ Row = 4;  E = -2029.447346335359; F = -302.9447346335359 },
    @{ Row = 5;  E = 1581.975128548893;  F = 58.19751285488933 },
    @{ Row = 7;  E = 1131.313681695051;  F = 13.1313681695051 },
    @{ Row = 8;  E = 1201.301932397386;  F = 20.13019323973865 },
    @{ Row = 9;  E = 871.0541417570648;  F = -12.89458582429351 },
    @{ Row = 10; E = 1058.680559134461;  F = 5.868055913446113 },
    @{ Row = 11; E = 1144.929498612026;  F = 14.49294986120265 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.E
    $ws.Cells.Item($u.Row, 6).Value = $u.F
}

$wb.Save()
